# TC47_Canine_Filter_Breed-YorkshireTerr.xlsx
# Fixed variables and query errors in Bread from TC30 to TC47
#
# The "CasesTab" row's Neo4j query (column B, row 2) incorrectly pulled in
# an extra `co:cohort` column that isn't part of this query's output. Strip
# the stray `Cohort` return column (and the OPTIONAL MATCH's now-unused
# trailing comma) so the query text matches its intended result set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

$fixedCaseQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Yorkshire Terrier']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $fixedCaseQuery

# Rows re-wrap/autofit to a shorter height now that the query text lost a line.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Selection moved off the old C4 cell onto the edited B2 cell, and the view
# no longer needs to be scrolled down to row 4.
[void]$ws.Range("B2").Select()
